{"js": "// Update the 25 \"two-digit \u00f7 one-digit\" answer cells in the single table.\n// Each populated table row (0, 4, 8, 12, 16) holds 5 answer cells (columns\n// 0-4); the rows in between are intentionally blank spacer rows. We replace\n// by explicit (row, col) position rather than text search, because a couple\n// of the old answer strings are duplicated elsewhere in the table but map to\n// DIFFERENT new values (e.g. \"63\u00f78=7, 7\" appears at (4,4) and (12,2), which\n// become \"42\u00f78=5, 2\" and \"86\u00f76=14, 2\" respectively).\nconst replacements = [\n  [0, 0, \"53\u00f75=10, 3\"],\n  [0, 1, \"16\u00f73=5, 1\"],\n  [0, 2, \"52\u00f73=17, 1\"],\n  [0, 3, \"10\u00f75=2, 0\"],\n  [0, 4, \"77\u00f76=12, 5\"],\n  [4, 0, \"79\u00f76=13, 1\"],\n  [4, 1, \"99\u00f77=14, 1\"],\n  [4, 2, \"83\u00f76=13, 5\"],\n  [4, 3, \"30\u00f79=3, 3\"],\n  [4, 4, \"42\u00f78=5, 2\"],\n  [8, 0, \"62\u00f79=6, 8\"],\n  [8, 1, \"27\u00f76=4, 3\"],\n  [8, 2, \"17\u00f73=5, 2\"],\n  [8, 3, \"56\u00f75=11, 1\"],\n  [8, 4, \"39\u00f75=7, 4\"],\n  [12, 0, \"61\u00f76=10, 1\"],\n  [12, 1, \"99\u00f74=24, 3\"],\n  [12, 2, \"86\u00f76=14, 2\"],\n  [12, 3, \"26\u00f75=5, 1\"],\n  [12, 4, \"62\u00f79=6, 8\"],\n  [16, 0, \"47\u00f79=5, 2\"],\n  [16, 1, \"68\u00f77=9, 5\"],\n  [16, 2, \"54\u00f73=18, 0\"],\n  [16, 3, \"77\u00f79=8, 5\"],\n  [16, 4, \"96\u00f74=24, 0\"],\n];\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\nfor (const [row, col, newText] of replacements) {\n  table.getCell(row, col).value = newText;\n}\nawait context.sync();\n", "ps1": "# Update the 25 \"two-digit \u00f7 one-digit\" answer cells in the single table.\n# Each populated table row (1, 5, 9, 13, 17 \u2014 1-based) holds 5 answer cells\n# (columns 1-5); the rows in between are intentionally blank spacer rows. We\n# replace by explicit (row, col) position rather than Find/Replace, because a\n# couple of the old answer strings are duplicated elsewhere in the table but\n# map to DIFFERENT new values (e.g. \"63\u00f78=7, 7\" appears at row5/col5 and\n# row13/col3, which become \"42\u00f78=5, 2\" and \"86\u00f76=14, 2\" respectively).\n\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n$replacements = @(\n  @{ Row = 1; Col = 1; Text = \"53\u00f75=10, 3\" },\n  @{ Row = 1; Col = 2; Text = \"16\u00f73=5, 1\" },\n  @{ Row = 1; Col = 3; Text = \"52\u00f73=17, 1\" },\n  @{ Row = 1; Col = 4; Text = \"10\u00f75=2, 0\" },\n  @{ Row = 1; Col = 5; Text = \"77\u00f76=12, 5\" },\n  @{ Row = 5; Col = 1; Text = \"79\u00f76=13, 1\" },\n  @{ Row = 5; Col = 2; Text = \"99\u00f77=14, 1\" },\n  @{ Row = 5; Col = 3; Text = \"83\u00f76=13, 5\" },\n  @{ Row = 5; Col = 4; Text = \"30\u00f79=3, 3\" },\n  @{ Row = 5; Col = 5; Text = \"42\u00f78=5, 2\" },\n  @{ Row = 9; Col = 1; Text = \"62\u00f79=6, 8\" },\n  @{ Row = 9; Col = 2; Text = \"27\u00f76=4, 3\" },\n  @{ Row = 9; Col = 3; Text = \"17\u00f73=5, 2\" },\n  @{ Row = 9; Col = 4; Text = \"56\u00f75=11, 1\" },\n  @{ Row = 9; Col = 5; Text = \"39\u00f75=7, 4\" },\n  @{ Row = 13; Col = 1; Text = \"61\u00f76=10, 1\" },\n  @{ Row = 13; Col = 2; Text = \"99\u00f74=24, 3\" },\n  @{ Row = 13; Col = 3; Text = \"86\u00f76=14, 2\" },\n  @{ Row = 13; Col = 4; Text = \"26\u00f75=5, 1\" },\n  @{ Row = 13; Col = 5; Text = \"62\u00f79=6, 8\" },\n  @{ Row = 17; Col = 1; Text = \"47\u00f79=5, 2\" },\n  @{ Row = 17; Col = 2; Text = \"68\u00f77=9, 5\" },\n  @{ Row = 17; Col = 3; Text = \"54\u00f73=18, 0\" },\n  @{ Row = 17; Col = 4; Text = \"77\u00f79=8, 5\" },\n  @{ Row = 17; Col = 5; Text = \"96\u00f74=24, 0\" }\n)\n\nforeach ($item in $replacements) {\n  $t.Cell($item.Row, $item.Col).Range.Text = $item.Text\n}\n"}
